$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the end of the data block (rows 97-98),
# shifting nothing existing (they are appended past the old last row).
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(98).Insert()

# Row 54
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = 'Vega Monumental Concepción'
$ws.Range("C54").Value = 'Bíobío'
$ws.Range("D54").Value = 44957
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = 'Frutos de hueso (carozo)'
$ws.Range("I54").Value = 100103002
$ws.Range("J54").Value = 'Ciruela'
$ws.Range("K54").Value = 'Fortuna'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 100
$ws.Range("N54").Value = 11000
$ws.Range("O54").Value = 12000
$ws.Range("P54").Value = 11500
$ws.Range("Q54").Value = '$/bandeja 18 kilos granel'
$ws.Range("R54").Value = 'Región de O''Higgins'
$ws.Range("S54").Value = 639
$ws.Range("T54").Value = 18

# Row 55
$ws.Range("A55").Value = 11
$ws.Range("B55").Value = 'Vega Monumental Concepción'
$ws.Range("C55").Value = 'Bíobío'
$ws.Range("D55").Value = 44957
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 'Fruta'
$ws.Range("G55").Value = 100103
$ws.Range("H55").Value = 'Frutos de hueso (carozo)'
$ws.Range("I55").Value = 100103002
$ws.Range("J55").Value = 'Ciruela'
$ws.Range("K55").Value = 'Fortuna'
$ws.Range("L55").Value = 'Segunda'
$ws.Range("M55").Value = 50
$ws.Range("N55").Value = 9000
$ws.Range("O55").Value = 9000
$ws.Range("P55").Value = 9000
$ws.Range("Q55").Value = '$/bandeja 18 kilos granel'
$ws.Range("R55").Value = 'Región de O''Higgins'
$ws.Range("S55").Value = 500
$ws.Range("T55").Value = 18

# Row 56
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = 'Vega Monumental Concepción'
$ws.Range("C56").Value = 'Bíobío'
$ws.Range("D56").Value = 44299
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 'Fruta'
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = 'Frutos de hueso (carozo)'
$ws.Range("I56").Value = 100103002
$ws.Range("J56").Value = 'Ciruela'
$ws.Range("K56").Value = 'Angeleno'
$ws.Range("L56").Value = 'Primera'
$ws.Range("M56").Value = 200
$ws.Range("N56").Value = 12000
$ws.Range("O56").Value = 13000
$ws.Range("P56").Value = 12500
$ws.Range("Q56").Value = '$/bandeja 18 kilos granel'
$ws.Range("R56").Value = 'Región de O''Higgins'
$ws.Range("S56").Value = 694
$ws.Range("T56").Value = 18

# Row 57
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = 'Vega Monumental Concepción'
$ws.Range("C57").Value = 'Bíobío'
$ws.Range("D57").Value = 44299
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 'Fruta'
$ws.Range("G57").Value = 100103
$ws.Range("H57").Value = 'Frutos de hueso (carozo)'
$ws.Range("I57").Value = 100103002
$ws.Range("J57").Value = 'Ciruela'
$ws.Range("K57").Value = 'Angeleno'
$ws.Range("L57").Value = 'Segunda'
$ws.Range("M57").Value = 100
$ws.Range("N57").Value = 11000
$ws.Range("O57").Value = 11000
$ws.Range("P57").Value = 11000
$ws.Range("Q57").Value = '$/bandeja 18 kilos granel'
$ws.Range("R57").Value = 'Región de O''Higgins'
$ws.Range("S57").Value = 611
$ws.Range("T57").Value = 18

# Row 58
$ws.Range("A58").Value = 11
$ws.Range("B58").Value = 'Vega Monumental Concepción'
$ws.Range("C58").Value = 'Bíobío'
$ws.Range("D58").Value = 44642
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 'Fruta'
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = 'Frutos de hueso (carozo)'
$ws.Range("I58").Value = 100103002
$ws.Range("J58").Value = 'Ciruela'
$ws.Range("K58").Value = 'Angeleno'
$ws.Range("L58").Value = 'Primera'
$ws.Range("M58").Value = 220
$ws.Range("N58").Value = 8000
$ws.Range("O58").Value = 9000
$ws.Range("P58").Value = 8455
$ws.Range("Q58").Value = '$/bandeja 18 kilos granel'
$ws.Range("R58").Value = 'Provincia de Curicó'
$ws.Range("S58").Value = 470
$ws.Range("T58").Value = 18

# Row 59
$ws.Range("A59").Value = 11
$ws.Range("B59").Value = 'Vega Monumental Concepción'
$ws.Range("C59").Value = 'Bíobío'
$ws.Range("D59").Value = 44588
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = 'Fruta'
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = 'Frutos de hueso (carozo)'
$ws.Range("I59").Value = 100103002
$ws.Range("J59").Value = 'Ciruela'
$ws.Range("K59").Value = 'Black Amber'
$ws.Range("L59").Value = 'Primera'
$ws.Range("M59").Value = 220
$ws.Range("N59").Value = 9000
$ws.Range("O59").Value = 9500
$ws.Range("P59").Value = 9227
$ws.Range("Q59").Value = '$/caja 16 kilos granel'
$ws.Range("R59").Value = 'Región de O''Higgins'
$ws.Range("S59").Value = 577
$ws.Range("T59").Value = 16

# Row 60
$ws.Range("A60").Value = 11
$ws.Range("B60").Value = 'Vega Monumental Concepción'
$ws.Range("C60").Value = 'Bíobío'
$ws.Range("D60").Value = 44202
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = 'Fruta'
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = 'Frutos de hueso (carozo)'
$ws.Range("I60").Value = 100103002
$ws.Range("J60").Value = 'Ciruela'
$ws.Range("K60").Value = 'Black Amber'
$ws.Range("L60").Value = 'Primera'
$ws.Range("M60").Value = 200
$ws.Range("N60").Value = 14000
$ws.Range("O60").Value = 15000
$ws.Range("P60").Value = 14500
$ws.Range("Q60").Value = '$/caja 18 kilos granel'
$ws.Range("R60").Value = 'Región de O''Higgins'
$ws.Range("S60").Value = 806
$ws.Range("T60").Value = 18

# Row 61
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = 'Vega Monumental Concepción'
$ws.Range("C61").Value = 'Bíobío'
$ws.Range("D61").Value = 44202
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 'Fruta'
$ws.Range("G61").Value = 100103
$ws.Range("H61").Value = 'Frutos de hueso (carozo)'
$ws.Range("I61").Value = 100103002
$ws.Range("J61").Value = 'Ciruela'
$ws.Range("K61").Value = 'Black Amber'
$ws.Range("L61").Value = 'Segunda'
$ws.Range("M61").Value = 100
$ws.Range("N61").Value = 12000
$ws.Range("O61").Value = 12000
$ws.Range("P61").Value = 12000
$ws.Range("Q61").Value = '$/caja 18 kilos granel'
$ws.Range("R61").Value = 'Región de O''Higgins'
$ws.Range("S61").Value = 667
$ws.Range("T61").Value = 18

# Row 62
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = 'Vega Monumental Concepción'
$ws.Range("C62").Value = 'Bíobío'
$ws.Range("D62").Value = 44931
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 'Fruta'
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = 'Frutos de hueso (carozo)'
$ws.Range("I62").Value = 100103002
$ws.Range("J62").Value = 'Ciruela'
$ws.Range("K62").Value = 'Black Amber'
$ws.Range("L62").Value = 'Primera'
$ws.Range("M62").Value = 100
$ws.Range("N62").Value = 15000
$ws.Range("O62").Value = 16000
$ws.Range("P62").Value = 15500
$ws.Range("Q62").Value = '$/bandeja 18 kilos granel'
$ws.Range("R62").Value = 'Región de O''Higgins'
$ws.Range("S62").Value = 861
$ws.Range("T62").Value = 18

# Row 63
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = 'Vega Monumental Concepción'
$ws.Range("C63").Value = 'Bíobío'
$ws.Range("D63").Value = 44931
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 'Fruta'
$ws.Range("G63").Value = 100103
$ws.Range("H63").Value = 'Frutos de hueso (carozo)'
$ws.Range("I63").Value = 100103002
$ws.Range("J63").Value = 'Ciruela'
$ws.Range("K63").Value = 'Black Amber'
$ws.Range("L63").Value = 'Segunda'
$ws.Range("M63").Value = 50
$ws.Range("N63").Value = 14000
$ws.Range("O63").Value = 14000
$ws.Range("P63").Value = 14000
$ws.Range("Q63").Value = '$/bandeja 18 kilos granel'
$ws.Range("R63").Value = 'Región de O''Higgins'
$ws.Range("S63").Value = 778
$ws.Range("T63").Value = 18

# Row 64
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = 'Vega Monumental Concepción'
$ws.Range("C64").Value = 'Bíobío'
$ws.Range("D64").Value = 44343
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = 'Fruta'
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = 'Frutos de hueso (carozo)'
$ws.Range("I64").Value = 100103002
$ws.Range("J64").Value = 'Ciruela'
$ws.Range("K64").Value = 'Angeleno'
$ws.Range("L64").Value = 'Primera'
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 11000
$ws.Range("P64").Value = 10500
$ws.Range("Q64").Value = '$/bandeja 18 kilos granel'
$ws.Range("R64").Value = 'Región de O''Higgins'
$ws.Range("S64").Value = 583
$ws.Range("T64").Value = 18

# Row 65
$ws.Range("A65").Value = 11
$ws.Range("B65").Value = 'Vega Monumental Concepción'
$ws.Range("C65").Value = 'Bíobío'
$ws.Range("D65").Value = 44343
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = 'Fruta'
$ws.Range("G65").Value = 100103
$ws.Range("H65").Value = 'Frutos de hueso (carozo)'
$ws.Range("I65").Value = 100103002
$ws.Range("J65").Value = 'Ciruela'
$ws.Range("K65").Value = 'Angeleno'
$ws.Range("L65").Value = 'Segunda'
$ws.Range("M65").Value = 100
$ws.Range("N65").Value = 9000
$ws.Range("O65").Value = 9000
$ws.Range("P65").Value = 9000
$ws.Range("Q65").Value = '$/bandeja 18 kilos granel'
$ws.Range("R65").Value = 'Región de O''Higgins'
$ws.Range("S65").Value = 500
$ws.Range("T65").Value = 18

# Row 66
$ws.Range("A66").Value = 11
$ws.Range("B66").Value = 'Vega Monumental Concepción'
$ws.Range("C66").Value = 'Bíobío'
$ws.Range("D66").Value = 44645
$ws.Range("E66").Value = 8
$ws.Range("F66").Value = 'Fruta'
$ws.Range("G66").Value = 100103
$ws.Range("H66").Value = 'Frutos de hueso (carozo)'
$ws.Range("I66").Value = 100103002
$ws.Range("J66").Value = 'Ciruela'
$ws.Range("K66").Value = 'Angeleno'
$ws.Range("L66").Value = 'Primera'
$ws.Range("M66").Value = 200
$ws.Range("N66").Value = 8000
$ws.Range("O66").Value = 8500
$ws.Range("P66").Value = 8250
$ws.Range("Q66").Value = '$/bandeja 18 kilos granel'
$ws.Range("R66").Value = 'Provincia de Curicó'
$ws.Range("S66").Value = 458
$ws.Range("T66").Value = 18

# Row 67
$ws.Range("A67").Value = 11
$ws.Range("B67").Value = 'Vega Monumental Concepción'
$ws.Range("C67").Value = 'Bíobío'
$ws.Range("D67").Value = 44645
$ws.Range("E67").Value = 8
$ws.Range("F67").Value = 'Fruta'
$ws.Range("G67").Value = 100103
$ws.Range("H67").Value = 'Frutos de hueso (carozo)'
$ws.Range("I67").Value = 100103002
$ws.Range("J67").Value = 'Ciruela'
$ws.Range("K67").Value = 'Angeleno'
$ws.Range("L67").Value = 'Segunda'
$ws.Range("M67").Value = 200
$ws.Range("N67").Value = 6000
$ws.Range("O67").Value = 6500
$ws.Range("P67").Value = 6250
$ws.Range("Q67").Value = '$/bandeja 18 kilos granel'
$ws.Range("R67").Value = 'Provincia de Curicó'
$ws.Range("S67").Value = 347
$ws.Range("T67").Value = 18

# Row 68
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = 'Vega Monumental Concepción'
$ws.Range("C68").Value = 'Bíobío'
$ws.Range("D68").Value = 44586
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = 'Fruta'
$ws.Range("G68").Value = 100103
$ws.Range("H68").Value = 'Frutos de hueso (carozo)'
$ws.Range("I68").Value = 100103002
$ws.Range("J68").Value = 'Ciruela'
$ws.Range("K68").Value = 'Black Amber'
$ws.Range("L68").Value = 'Primera'
$ws.Range("M68").Value = 200
$ws.Range("N68").Value = 9000
$ws.Range("O68").Value = 10000
$ws.Range("P68").Value = 9500
$ws.Range("Q68").Value = '$/bandeja 18 kilos granel'
$ws.Range("R68").Value = 'Región de O''Higgins'
$ws.Range("S68").Value = 528
$ws.Range("T68").Value = 18

# Row 69
$ws.Range("A69").Value = 11
$ws.Range("B69").Value = 'Vega Monumental Concepción'
$ws.Range("C69").Value = 'Bíobío'
$ws.Range("D69").Value = 44586
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 'Fruta'
$ws.Range("G69").Value = 100103
$ws.Range("H69").Value = 'Frutos de hueso (carozo)'
$ws.Range("I69").Value = 100103002
$ws.Range("J69").Value = 'Ciruela'
$ws.Range("K69").Value = 'Black Amber'
$ws.Range("L69").Value = 'Segunda'
$ws.Range("M69").Value = 100
$ws.Range("N69").Value = 8000
$ws.Range("O69").Value = 8000
$ws.Range("P69").Value = 8000
$ws.Range("Q69").Value = '$/bandeja 18 kilos granel'
$ws.Range("R69").Value = 'Región de O''Higgins'
$ws.Range("S69").Value = 444
$ws.Range("T69").Value = 18

# Row 70
$ws.Range("A70").Value = 11
$ws.Range("B70").Value = 'Vega Monumental Concepción'
$ws.Range("C70").Value = 'Bíobío'
$ws.Range("D70").Value = 44307
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = 'Fruta'
$ws.Range("G70").Value = 100103
$ws.Range("H70").Value = 'Frutos de hueso (carozo)'
$ws.Range("I70").Value = 100103002
$ws.Range("J70").Value = 'Ciruela'
$ws.Range("K70").Value = 'Angeleno'
$ws.Range("L70").Value = 'Primera'
$ws.Range("M70").Value = 200
$ws.Range("N70").Value = 9000
$ws.Range("O70").Value = 10000
$ws.Range("P70").Value = 9500
$ws.Range("Q70").Value = '$/bandeja 18 kilos granel'
$ws.Range("R70").Value = 'Región de O''Higgins'
$ws.Range("S70").Value = 528
$ws.Range("T70").Value = 18

# Row 71
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = 'Vega Monumental Concepción'
$ws.Range("C71").Value = 'Bíobío'
$ws.Range("D71").Value = 44307
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 'Fruta'
$ws.Range("G71").Value = 100103
$ws.Range("H71").Value = 'Frutos de hueso (carozo)'
$ws.Range("I71").Value = 100103002
$ws.Range("J71").Value = 'Ciruela'
$ws.Range("K71").Value = 'Angeleno'
$ws.Range("L71").Value = 'Segunda'
$ws.Range("M71").Value = 100
$ws.Range("N71").Value = 8000
$ws.Range("O71").Value = 8000
$ws.Range("P71").Value = 8000
$ws.Range("Q71").Value = '$/bandeja 18 kilos granel'
$ws.Range("R71").Value = 'Región de O''Higgins'
$ws.Range("S71").Value = 444
$ws.Range("T71").Value = 18

# Row 72
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = 'Vega Monumental Concepción'
$ws.Range("C72").Value = 'Bíobío'
$ws.Range("D72").Value = 44215
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 'Fruta'
$ws.Range("G72").Value = 100103
$ws.Range("H72").Value = 'Frutos de hueso (carozo)'
$ws.Range("I72").Value = 100103002
$ws.Range("J72").Value = 'Ciruela'
$ws.Range("K72").Value = 'Black Amber'
$ws.Range("L72").Value = 'Primera'
$ws.Range("M72").Value = 200
$ws.Range("N72").Value = 10000
$ws.Range("O72").Value = 11000
$ws.Range("P72").Value = 10500
$ws.Range("Q72").Value = '$/caja 16 kilos granel'
$ws.Range("R72").Value = 'Región de O''Higgins'
$ws.Range("S72").Value = 656
$ws.Range("T72").Value = 16

# Row 73
$ws.Range("A73").Value = 11
$ws.Range("B73").Value = 'Vega Monumental Concepción'
$ws.Range("C73").Value = 'Bíobío'
$ws.Range("D73").Value = 44215
$ws.Range("E73").Value = 8
$ws.Range("F73").Value = 'Fruta'
$ws.Range("G73").Value = 100103
$ws.Range("H73").Value = 'Frutos de hueso (carozo)'
$ws.Range("I73").Value = 100103002
$ws.Range("J73").Value = 'Ciruela'
$ws.Range("K73").Value = 'Black Amber'
$ws.Range("L73").Value = 'Segunda'
$ws.Range("M73").Value = 100
$ws.Range("N73").Value = 8000
$ws.Range("O73").Value = 8000
$ws.Range("P73").Value = 8000
$ws.Range("Q73").Value = '$/caja 16 kilos granel'
$ws.Range("R73").Value = 'Región de O''Higgins'
$ws.Range("S73").Value = 500
$ws.Range("T73").Value = 16

# Row 74
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = 'Vega Monumental Concepción'
$ws.Range("C74").Value = 'Bíobío'
$ws.Range("D74").Value = 44951
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 'Fruta'
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = 'Frutos de hueso (carozo)'
$ws.Range("I74").Value = 100103002
$ws.Range("J74").Value = 'Ciruela'
$ws.Range("K74").Value = 'Fortuna'
$ws.Range("L74").Value = 'Primera'
$ws.Range("M74").Value = 170
$ws.Range("N74").Value = 10000
$ws.Range("O74").Value = 11000
$ws.Range("P74").Value = 10529
$ws.Range("Q74").Value = '$/caja 18 kilos granel'
$ws.Range("R74").Value = 'Provincia de Curicó'
$ws.Range("S74").Value = 585
$ws.Range("T74").Value = 18

# Row 75
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = 'Vega Monumental Concepción'
$ws.Range("C75").Value = 'Bíobío'
$ws.Range("D75").Value = 44266
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 'Fruta'
$ws.Range("G75").Value = 100103
$ws.Range("H75").Value = 'Frutos de hueso (carozo)'
$ws.Range("I75").Value = 100103002
$ws.Range("J75").Value = 'Ciruela'
$ws.Range("K75").Value = 'Black Amber'
$ws.Range("L75").Value = 'Primera'
$ws.Range("M75").Value = 200
$ws.Range("N75").Value = 9000
$ws.Range("O75").Value = 10000
$ws.Range("P75").Value = 9500
$ws.Range("Q75").Value = '$/caja 18 kilos granel'
$ws.Range("R75").Value = 'Región de O''Higgins'
$ws.Range("S75").Value = 528
$ws.Range("T75").Value = 18

# Row 76
$ws.Range("A76").Value = 11
$ws.Range("B76").Value = 'Vega Monumental Concepción'
$ws.Range("C76").Value = 'Bíobío'
$ws.Range("D76").Value = 44266
$ws.Range("E76").Value = 8
$ws.Range("F76").Value = 'Fruta'
$ws.Range("G76").Value = 100103
$ws.Range("H76").Value = 'Frutos de hueso (carozo)'
$ws.Range("I76").Value = 100103002
$ws.Range("J76").Value = 'Ciruela'
$ws.Range("K76").Value = 'Black Amber'
$ws.Range("L76").Value = 'Segunda'
$ws.Range("M76").Value = 100
$ws.Range("N76").Value = 8000
$ws.Range("O76").Value = 8000
$ws.Range("P76").Value = 8000
$ws.Range("Q76").Value = '$/caja 18 kilos granel'
$ws.Range("R76").Value = 'Región de O''Higgins'
$ws.Range("S76").Value = 444
$ws.Range("T76").Value = 18

# Row 77
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = 'Vega Monumental Concepción'
$ws.Range("C77").Value = 'Bíobío'
$ws.Range("D77").Value = 44607
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = 'Fruta'
$ws.Range("G77").Value = 100103
$ws.Range("H77").Value = 'Frutos de hueso (carozo)'
$ws.Range("I77").Value = 100103002
$ws.Range("J77").Value = 'Ciruela'
$ws.Range("K77").Value = 'Black Amber'
$ws.Range("L77").Value = 'Primera'
$ws.Range("M77").Value = 250
$ws.Range("N77").Value = 11000
$ws.Range("O77").Value = 12000
$ws.Range("P77").Value = 11520
$ws.Range("Q77").Value = '$/bandeja 18 kilos granel'
$ws.Range("R77").Value = 'Región de O''Higgins'
$ws.Range("S77").Value = 640
$ws.Range("T77").Value = 18

# Row 78
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = 'Vega Monumental Concepción'
$ws.Range("C78").Value = 'Bíobío'
$ws.Range("D78").Value = 44607
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 'Fruta'
$ws.Range("G78").Value = 100103
$ws.Range("H78").Value = 'Frutos de hueso (carozo)'
$ws.Range("I78").Value = 100103002
$ws.Range("J78").Value = 'Ciruela'
$ws.Range("K78").Value = 'Black Amber'
$ws.Range("L78").Value = 'Segunda'
$ws.Range("M78").Value = 170
$ws.Range("N78").Value = 9000
$ws.Range("O78").Value = 9500
$ws.Range("P78").Value = 9265
$ws.Range("Q78").Value = '$/bandeja 18 kilos granel'
$ws.Range("R78").Value = 'Región de O''Higgins'
$ws.Range("S78").Value = 515
$ws.Range("T78").Value = 18

# Row 79
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = 'Vega Monumental Concepción'
$ws.Range("C79").Value = 'Bíobío'
$ws.Range("D79").Value = 44328
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = 'Fruta'
$ws.Range("G79").Value = 100103
$ws.Range("H79").Value = 'Frutos de hueso (carozo)'
$ws.Range("I79").Value = 100103002
$ws.Range("J79").Value = 'Ciruela'
$ws.Range("K79").Value = 'Angeleno'
$ws.Range("L79").Value = 'Primera'
$ws.Range("M79").Value = 100
$ws.Range("N79").Value = 9000
$ws.Range("O79").Value = 10000
$ws.Range("P79").Value = 9500
$ws.Range("Q79").Value = '$/bandeja 18 kilos granel'
$ws.Range("R79").Value = 'Región de O''Higgins'
$ws.Range("S79").Value = 528
$ws.Range("T79").Value = 18

# Row 80
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = 'Vega Monumental Concepción'
$ws.Range("C80").Value = 'Bíobío'
$ws.Range("D80").Value = 44328
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = 'Fruta'
$ws.Range("G80").Value = 100103
$ws.Range("H80").Value = 'Frutos de hueso (carozo)'
$ws.Range("I80").Value = 100103002
$ws.Range("J80").Value = 'Ciruela'
$ws.Range("K80").Value = 'Angeleno'
$ws.Range("L80").Value = 'Segunda'
$ws.Range("M80").Value = 50
$ws.Range("N80").Value = 8000
$ws.Range("O80").Value = 8000
$ws.Range("P80").Value = 8000
$ws.Range("Q80").Value = '$/bandeja 18 kilos granel'
$ws.Range("R80").Value = 'Región de O''Higgins'
$ws.Range("S80").Value = 444
$ws.Range("T80").Value = 18

# Row 81
$ws.Range("A81").Value = 11
$ws.Range("B81").Value = 'Vega Monumental Concepción'
$ws.Range("C81").Value = 'Bíobío'
$ws.Range("D81").Value = 44946
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 'Fruta'
$ws.Range("G81").Value = 100103
$ws.Range("H81").Value = 'Frutos de hueso (carozo)'
$ws.Range("I81").Value = 100103002
$ws.Range("J81").Value = 'Ciruela'
$ws.Range("K81").Value = 'Black Amber'
$ws.Range("L81").Value = 'Primera'
$ws.Range("M81").Value = 100
$ws.Range("N81").Value = 11000
$ws.Range("O81").Value = 12000
$ws.Range("P81").Value = 11500
$ws.Range("Q81").Value = '$/bandeja 18 kilos granel'
$ws.Range("R81").Value = 'Región de O''Higgins'
$ws.Range("S81").Value = 639
$ws.Range("T81").Value = 18

# Row 82
$ws.Range("A82").Value = 11
$ws.Range("B82").Value = 'Vega Monumental Concepción'
$ws.Range("C82").Value = 'Bíobío'
$ws.Range("D82").Value = 44946
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = 'Fruta'
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = 'Frutos de hueso (carozo)'
$ws.Range("I82").Value = 100103002
$ws.Range("J82").Value = 'Ciruela'
$ws.Range("K82").Value = 'Black Amber'
$ws.Range("L82").Value = 'Segunda'
$ws.Range("M82").Value = 50
$ws.Range("N82").Value = 9000
$ws.Range("O82").Value = 9000
$ws.Range("P82").Value = 9000
$ws.Range("Q82").Value = '$/bandeja 18 kilos granel'
$ws.Range("R82").Value = 'Región de O''Higgins'
$ws.Range("S82").Value = 500
$ws.Range("T82").Value = 18

# Row 83
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = 'Vega Monumental Concepción'
$ws.Range("C83").Value = 'Bíobío'
$ws.Range("D83").Value = 44644
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 'Fruta'
$ws.Range("G83").Value = 100103
$ws.Range("H83").Value = 'Frutos de hueso (carozo)'
$ws.Range("I83").Value = 100103002
$ws.Range("J83").Value = 'Ciruela'
$ws.Range("K83").Value = 'Angeleno'
$ws.Range("L83").Value = 'Primera'
$ws.Range("M83").Value = 350
$ws.Range("N83").Value = 8500
$ws.Range("O83").Value = 9000
$ws.Range("P83").Value = 8714
$ws.Range("Q83").Value = '$/bandeja 18 kilos granel'
$ws.Range("R83").Value = 'Provincia de Curicó'
$ws.Range("S83").Value = 484
$ws.Range("T83").Value = 18

# Row 84
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = 'Vega Monumental Concepción'
$ws.Range("C84").Value = 'Bíobío'
$ws.Range("D84").Value = 44637
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 'Fruta'
$ws.Range("G84").Value = 100103
$ws.Range("H84").Value = 'Frutos de hueso (carozo)'
$ws.Range("I84").Value = 100103002
$ws.Range("J84").Value = 'Ciruela'
$ws.Range("K84").Value = 'Angeleno'
$ws.Range("L84").Value = 'Primera'
$ws.Range("M84").Value = 220
$ws.Range("N84").Value = 8000
$ws.Range("O84").Value = 8500
$ws.Range("P84").Value = 8273
$ws.Range("Q84").Value = '$/bandeja 18 kilos granel'
$ws.Range("R84").Value = 'Provincia de Curicó'
$ws.Range("S84").Value = 460
$ws.Range("T84").Value = 18

# Row 85
$ws.Range("A85").Value = 11
$ws.Range("B85").Value = 'Vega Monumental Concepción'
$ws.Range("C85").Value = 'Bíobío'
$ws.Range("D85").Value = 44223
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = 'Fruta'
$ws.Range("G85").Value = 100103
$ws.Range("H85").Value = 'Frutos de hueso (carozo)'
$ws.Range("I85").Value = 100103002
$ws.Range("J85").Value = 'Ciruela'
$ws.Range("K85").Value = 'Black Amber'
$ws.Range("L85").Value = 'Primera'
$ws.Range("M85").Value = 100
$ws.Range("N85").Value = 10000
$ws.Range("O85").Value = 11000
$ws.Range("P85").Value = 10500
$ws.Range("Q85").Value = '$/caja 16 kilos granel'
$ws.Range("R85").Value = 'Región de O''Higgins'
$ws.Range("S85").Value = 656
$ws.Range("T85").Value = 16

# Row 86
$ws.Range("A86").Value = 11
$ws.Range("B86").Value = 'Vega Monumental Concepción'
$ws.Range("C86").Value = 'Bíobío'
$ws.Range("D86").Value = 44223
$ws.Range("E86").Value = 8
$ws.Range("F86").Value = 'Fruta'
$ws.Range("G86").Value = 100103
$ws.Range("H86").Value = 'Frutos de hueso (carozo)'
$ws.Range("I86").Value = 100103002
$ws.Range("J86").Value = 'Ciruela'
$ws.Range("K86").Value = 'Black Amber'
$ws.Range("L86").Value = 'Segunda'
$ws.Range("M86").Value = 50
$ws.Range("N86").Value = 9000
$ws.Range("O86").Value = 9000
$ws.Range("P86").Value = 9000
$ws.Range("Q86").Value = '$/caja 16 kilos granel'
$ws.Range("R86").Value = 'Región de O''Higgins'
$ws.Range("S86").Value = 562
$ws.Range("T86").Value = 16

# Row 87
$ws.Range("A87").Value = 11
$ws.Range("B87").Value = 'Vega Monumental Concepción'
$ws.Range("C87").Value = 'Bíobío'
$ws.Range("D87").Value = 44616
$ws.Range("E87").Value = 8
$ws.Range("F87").Value = 'Fruta'
$ws.Range("G87").Value = 100103
$ws.Range("H87").Value = 'Frutos de hueso (carozo)'
$ws.Range("I87").Value = 100103002
$ws.Range("J87").Value = 'Ciruela'
$ws.Range("K87").Value = 'Black Amber'
$ws.Range("L87").Value = 'Primera'
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = 9000
$ws.Range("O87").Value = 10000
$ws.Range("P87").Value = 9500
$ws.Range("Q87").Value = '$/bandeja 18 kilos granel'
$ws.Range("R87").Value = 'Región de O''Higgins'
$ws.Range("S87").Value = 528
$ws.Range("T87").Value = 18

# Row 88
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = 'Vega Monumental Concepción'
$ws.Range("C88").Value = 'Bíobío'
$ws.Range("D88").Value = 44616
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 'Fruta'
$ws.Range("G88").Value = 100103
$ws.Range("H88").Value = 'Frutos de hueso (carozo)'
$ws.Range("I88").Value = 100103002
$ws.Range("J88").Value = 'Ciruela'
$ws.Range("K88").Value = 'Black Amber'
$ws.Range("L88").Value = 'Segunda'
$ws.Range("M88").Value = 100
$ws.Range("N88").Value = 8000
$ws.Range("O88").Value = 8000
$ws.Range("P88").Value = 8000
$ws.Range("Q88").Value = '$/bandeja 18 kilos granel'
$ws.Range("R88").Value = 'Región de O''Higgins'
$ws.Range("S88").Value = 444
$ws.Range("T88").Value = 18

# Row 89
$ws.Range("A89").Value = 11
$ws.Range("B89").Value = 'Vega Monumental Concepción'
$ws.Range("C89").Value = 'Bíobío'
$ws.Range("D89").Value = 44631
$ws.Range("E89").Value = 8
$ws.Range("F89").Value = 'Fruta'
$ws.Range("G89").Value = 100103
$ws.Range("H89").Value = 'Frutos de hueso (carozo)'
$ws.Range("I89").Value = 100103002
$ws.Range("J89").Value = 'Ciruela'
$ws.Range("K89").Value = 'Angeleno'
$ws.Range("L89").Value = 'Primera'
$ws.Range("M89").Value = 100
$ws.Range("N89").Value = 8000
$ws.Range("O89").Value = 9000
$ws.Range("P89").Value = 8500
$ws.Range("Q89").Value = '$/bandeja 18 kilos granel'
$ws.Range("R89").Value = 'Región de O''Higgins'
$ws.Range("S89").Value = 472
$ws.Range("T89").Value = 18

# Row 90
$ws.Range("A90").Value = 11
$ws.Range("B90").Value = 'Vega Monumental Concepción'
$ws.Range("C90").Value = 'Bíobío'
$ws.Range("D90").Value = 44631
$ws.Range("E90").Value = 8
$ws.Range("F90").Value = 'Fruta'
$ws.Range("G90").Value = 100103
$ws.Range("H90").Value = 'Frutos de hueso (carozo)'
$ws.Range("I90").Value = 100103002
$ws.Range("J90").Value = 'Ciruela'
$ws.Range("K90").Value = 'Angeleno'
$ws.Range("L90").Value = 'Segunda'
$ws.Range("M90").Value = 50
$ws.Range("N90").Value = 7000
$ws.Range("O90").Value = 7000
$ws.Range("P90").Value = 7000
$ws.Range("Q90").Value = '$/bandeja 18 kilos granel'
$ws.Range("R90").Value = 'Región de O''Higgins'
$ws.Range("S90").Value = 389
$ws.Range("T90").Value = 18

# Row 91
$ws.Range("A91").Value = 11
$ws.Range("B91").Value = 'Vega Monumental Concepción'
$ws.Range("C91").Value = 'Bíobío'
$ws.Range("D91").Value = 44643
$ws.Range("E91").Value = 8
$ws.Range("F91").Value = 'Fruta'
$ws.Range("G91").Value = 100103
$ws.Range("H91").Value = 'Frutos de hueso (carozo)'
$ws.Range("I91").Value = 100103002
$ws.Range("J91").Value = 'Ciruela'
$ws.Range("K91").Value = 'Angeleno'
$ws.Range("L91").Value = 'Primera'
$ws.Range("M91").Value = 180
$ws.Range("N91").Value = 8500
$ws.Range("O91").Value = 9000
$ws.Range("P91").Value = 8778
$ws.Range("Q91").Value = '$/bandeja 18 kilos granel'
$ws.Range("R91").Value = 'Región del Maule'
$ws.Range("S91").Value = 488
$ws.Range("T91").Value = 18

# Row 92
$ws.Range("A92").Value = 11
$ws.Range("B92").Value = 'Vega Monumental Concepción'
$ws.Range("C92").Value = 'Bíobío'
$ws.Range("D92").Value = 44643
$ws.Range("E92").Value = 8
$ws.Range("F92").Value = 'Fruta'
$ws.Range("G92").Value = 100103
$ws.Range("H92").Value = 'Frutos de hueso (carozo)'
$ws.Range("I92").Value = 100103002
$ws.Range("J92").Value = 'Ciruela'
$ws.Range("K92").Value = 'Angeleno'
$ws.Range("L92").Value = 'Segunda'
$ws.Range("M92").Value = 150
$ws.Range("N92").Value = 6500
$ws.Range("O92").Value = 7000
$ws.Range("P92").Value = 6767
$ws.Range("Q92").Value = '$/bandeja 18 kilos granel'
$ws.Range("R92").Value = 'Región del Maule'
$ws.Range("S92").Value = 376
$ws.Range("T92").Value = 18

# Row 93
$ws.Range("A93").Value = 11
$ws.Range("B93").Value = 'Vega Monumental Concepción'
$ws.Range("C93").Value = 'Bíobío'
$ws.Range("D93").Value = 44657
$ws.Range("E93").Value = 8
$ws.Range("F93").Value = 'Fruta'
$ws.Range("G93").Value = 100103
$ws.Range("H93").Value = 'Frutos de hueso (carozo)'
$ws.Range("I93").Value = 100103002
$ws.Range("J93").Value = 'Ciruela'
$ws.Range("K93").Value = 'Angeleno'
$ws.Range("L93").Value = 'Primera'
$ws.Range("M93").Value = 200
$ws.Range("N93").Value = 9000
$ws.Range("O93").Value = 10000
$ws.Range("P93").Value = 9500
$ws.Range("Q93").Value = '$/bandeja 18 kilos granel'
$ws.Range("R93").Value = 'Región de O''Higgins'
$ws.Range("S93").Value = 528
$ws.Range("T93").Value = 18

# Row 94
$ws.Range("A94").Value = 11
$ws.Range("B94").Value = 'Vega Monumental Concepción'
$ws.Range("C94").Value = 'Bíobío'
$ws.Range("D94").Value = 44657
$ws.Range("E94").Value = 8
$ws.Range("F94").Value = 'Fruta'
$ws.Range("G94").Value = 100103
$ws.Range("H94").Value = 'Frutos de hueso (carozo)'
$ws.Range("I94").Value = 100103002
$ws.Range("J94").Value = 'Ciruela'
$ws.Range("K94").Value = 'Angeleno'
$ws.Range("L94").Value = 'Segunda'
$ws.Range("M94").Value = 100
$ws.Range("N94").Value = 8000
$ws.Range("O94").Value = 8000
$ws.Range("P94").Value = 8000
$ws.Range("Q94").Value = '$/bandeja 18 kilos granel'
$ws.Range("R94").Value = 'Región de O''Higgins'
$ws.Range("S94").Value = 444
$ws.Range("T94").Value = 18

# Row 95
$ws.Range("A95").Value = 11
$ws.Range("B95").Value = 'Vega Monumental Concepción'
$ws.Range("C95").Value = 'Bíobío'
$ws.Range("D95").Value = 44251
$ws.Range("E95").Value = 8
$ws.Range("F95").Value = 'Fruta'
$ws.Range("G95").Value = 100103
$ws.Range("H95").Value = 'Frutos de hueso (carozo)'
$ws.Range("I95").Value = 100103002
$ws.Range("J95").Value = 'Ciruela'
$ws.Range("K95").Value = 'Angeleno'
$ws.Range("L95").Value = 'Primera'
$ws.Range("M95").Value = 200
$ws.Range("N95").Value = 9000
$ws.Range("O95").Value = 10000
$ws.Range("P95").Value = 9500
$ws.Range("Q95").Value = '$/caja 16 kilos granel'
$ws.Range("R95").Value = 'Región de O''Higgins'
$ws.Range("S95").Value = 594
$ws.Range("T95").Value = 16

# Row 96
$ws.Range("A96").Value = 11
$ws.Range("B96").Value = 'Vega Monumental Concepción'
$ws.Range("C96").Value = 'Bíobío'
$ws.Range("D96").Value = 44251
$ws.Range("E96").Value = 8
$ws.Range("F96").Value = 'Fruta'
$ws.Range("G96").Value = 100103
$ws.Range("H96").Value = 'Frutos de hueso (carozo)'
$ws.Range("I96").Value = 100103002
$ws.Range("J96").Value = 'Ciruela'
$ws.Range("K96").Value = 'Angeleno'
$ws.Range("L96").Value = 'Segunda'
$ws.Range("M96").Value = 100
$ws.Range("N96").Value = 8000
$ws.Range("O96").Value = 8000
$ws.Range("P96").Value = 8000
$ws.Range("Q96").Value = '$/caja 16 kilos granel'
$ws.Range("R96").Value = 'Región de O''Higgins'
$ws.Range("S96").Value = 500
$ws.Range("T96").Value = 16

# Row 97
$ws.Range("A97").Value = 11
$ws.Range("B97").Value = 'Vega Monumental Concepción'
$ws.Range("C97").Value = 'Bíobío'
$ws.Range("D97").Value = 44636
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 'Fruta'
$ws.Range("G97").Value = 100103
$ws.Range("H97").Value = 'Frutos de hueso (carozo)'
$ws.Range("I97").Value = 100103002
$ws.Range("J97").Value = 'Ciruela'
$ws.Range("K97").Value = 'Angeleno'
$ws.Range("L97").Value = 'Segunda'
$ws.Range("M97").Value = 220
$ws.Range("N97").Value = 6500
$ws.Range("O97").Value = 7000
$ws.Range("P97").Value = 6727
$ws.Range("Q97").Value = '$/bandeja 18 kilos granel'
$ws.Range("R97").Value = 'Provincia de Curicó'
$ws.Range("S97").Value = 374
$ws.Range("T97").Value = 18

# Row 98
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = 'Vega Monumental Concepción'
$ws.Range("C98").Value = 'Bíobío'
$ws.Range("D98").Value = 44595
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 'Fruta'
$ws.Range("G98").Value = 100103
$ws.Range("H98").Value = 'Frutos de hueso (carozo)'
$ws.Range("I98").Value = 100103002
$ws.Range("J98").Value = 'Ciruela'
$ws.Range("K98").Value = 'Black Amber'
$ws.Range("L98").Value = 'Primera'
$ws.Range("M98").Value = 250
$ws.Range("N98").Value = 8500
$ws.Range("O98").Value = 9000
$ws.Range("P98").Value = 8740
$ws.Range("Q98").Value = '$/caja 16 kilos granel'
$ws.Range("R98").Value = 'Provincia de Curicó'
$ws.Range("S98").Value = 546
$ws.Range("T98").Value = 16
